# "Generate Report for Handback" — populate the Latest Target File / Latest
# Handback File / Latest Handback DateTime columns (F/G/H) for the zh-cn and
# de-de detail sheets, and flip the Status column from "Ready for handoff" to
# "Handed back: in sync with en-US" everywhere it appears (Overview + both
# detail sheets).

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# Overview sheet: Status columns B (zh-cn) and C (de-de) for both file rows.
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $newStatus
$overview.Range("C2").Value = $newStatus
$overview.Range("B3").Value = $newStatus
$overview.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------------
# Per-locale detail sheets. Each row's "Latest Target File" (F) mirrors the
# "Source File Name" (A) hyperlink/text, and "Latest Handback File" (G)
# mirrors the "Latest Handoff File" (D) hyperlink/text — the handback just
# shipped back the same target file that was originally handed off.
# ---------------------------------------------------------------------------
$rows = @(2, 3)

$zhCn = $wb.Worksheets.Item("zh-cn")
$zhCnHandbackTime = "2016-03-18 05:25:20"
$zhCnTargetUrls = @{
    2 = "https://github.com/OpenLocalizationTest/oltest/blob/a5fd9243243b4a8b6d64e5e34cce16a8c5777b0b/e2e/2eceaf31-de09-4ec1-8476-9046d604730a.md"
    3 = "https://github.com/OpenLocalizationTest/oltest/blob/a5fd9243243b4a8b6d64e5e34cce16a8c5777b0b/e2e/785d7a05-0861-4070-b95a-8f1a31d8450b.md"
}
$zhCnHandbackUrls = @{
    2 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e0fc7474f9a6268c47f4a7ec36e8238803c9883b/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/high/2eceaf31-de09-4ec1-8476-9046d604730a.101ca9cd84e24e398a9261d96faf06407c4f9ba7.zh-cn.xlf"
    3 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e0fc7474f9a6268c47f4a7ec36e8238803c9883b/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/high/785d7a05-0861-4070-b95a-8f1a31d8450b.8dcbc33b3cec6fb1f7c2ce38db43d91d78d9d6ee.zh-cn.xlf"
}

$deDe = $wb.Worksheets.Item("de-de")
$deDeHandbackTime = "2016-03-18 05:25:25"
$deDeTargetUrls = @{
    2 = "https://github.com/OpenLocalizationTest/oltest/blob/a5fd9243243b4a8b6d64e5e34cce16a8c5777b0b/e2e/2eceaf31-de09-4ec1-8476-9046d604730a.md"
    3 = "https://github.com/OpenLocalizationTest/oltest/blob/a5fd9243243b4a8b6d64e5e34cce16a8c5777b0b/e2e/785d7a05-0861-4070-b95a-8f1a31d8450b.md"
}
$deDeHandbackUrls = @{
    2 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/772e11b502a7df2f302c72864d2ff4fecd8fc363/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/high/2eceaf31-de09-4ec1-8476-9046d604730a.101ca9cd84e24e398a9261d96faf06407c4f9ba7.de-de.xlf"
    3 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/772e11b502a7df2f302c72864d2ff4fecd8fc363/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/high/785d7a05-0861-4070-b95a-8f1a31d8450b.8dcbc33b3cec6fb1f7c2ce38db43d91d78d9d6ee.de-de.xlf"
}

function Fill-HandbackRow {
    param(
        $ws,
        [int]$row,
        [string]$handbackTime,
        [string]$targetUrl,
        [string]$handbackUrl
    )

    # Status
    $ws.Range("C$row").Value = $newStatus

    # Latest Target File (F) — same display text + link as Source File Name (A)
    $srcText = $ws.Range("A$row").Text
    $fCell = $ws.Range("F$row")
    $fCell.Value = $srcText
    $ws.Hyperlinks.Add($fCell, $targetUrl, [Type]::Missing, [Type]::Missing, $srcText) | Out-Null
    $fCell.Style = "HyperLink"

    # Latest Handback File (G) — same display text + link as Latest Handoff File (D)
    $handoffText = $ws.Range("D$row").Text
    $gCell = $ws.Range("G$row")
    $gCell.Value = $handoffText
    $ws.Hyperlinks.Add($gCell, $handbackUrl, [Type]::Missing, [Type]::Missing, $handoffText) | Out-Null
    $gCell.Style = "HyperLink"

    # Latest Handback DateTime (H)
    $ws.Range("H$row").Value = $handbackTime
}

foreach ($row in $rows) {
    Fill-HandbackRow $zhCn $row $zhCnHandbackTime $zhCnTargetUrls[$row] $zhCnHandbackUrls[$row]
    Fill-HandbackRow $deDe $row $deDeHandbackTime $deDeTargetUrls[$row] $deDeHandbackUrls[$row]
}

Write-Host "Handback report generated."
